# Append one new data row (row 67) to each of the four sheets, mirroring the
# existing row-67 records added to the logging database export.
# Column layout: A=time, B=总长, C=ID, D=实际长度, E=和校验,
#                F=总长_DEC, G=ID_DEC, H=实际长度_DEC, I=和校验_DEC

$wb = $excel.ActiveWorkbook

$rows = @{
    "DE_LFT_#1" = @{
        A = 45853.43400462963
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x50"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 336
        I = 14
    }
    "DE_LFT_#2" = @{
        A = 45853.43400462963
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x50"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 336
        I = 14
    }
    "DE_PLT_#1" = @{
        A = 45853.43400462963
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7B"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 123
        I = 7
    }
    "DE_PLT_#2" = @{
        A = 45853.43400462963
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7B"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 123
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rows[$ws.Name]
    if ($data -eq $null) { continue }

    $newRow = 67

    # New row's date cell should carry the same date/time number format as
    # the cell directly above it.
    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
